# Insert a new weekly record at row 6.
# This pushes the existing rows 6..66 down to 7..67 (preserving all of
# their data/formatting untouched) and fills the freshly inserted row 6
# with the new observation. Several columns (A,B,C,E,F,G,H,I,J,K,Q,R,T)
# keep the same values the old row 6 had, so we simply read them back
# from row 7 (which now holds what used to be row 6) after the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (6-66) down by one, inserting a blank row at 6.
$ws.Rows.Item(6).Insert()

# Columns that stay the same as the row that used to occupy position 6
# (now shifted to row 7).
$ws.Cells.Item(6, 1).Value = $ws.Cells.Item(7, 1).Value2   # Mercado ID
$ws.Cells.Item(6, 2).Value = $ws.Cells.Item(7, 2).Value2   # Mercado
$ws.Cells.Item(6, 3).Value = $ws.Cells.Item(7, 3).Value2   # Región
$ws.Cells.Item(6, 5).Value = $ws.Cells.Item(7, 5).Value2   # Codreg
$ws.Cells.Item(6, 6).Value = $ws.Cells.Item(7, 6).Value2   # Tipo
$ws.Cells.Item(6, 7).Value = $ws.Cells.Item(7, 7).Value2   # Producto ID
$ws.Cells.Item(6, 8).Value = $ws.Cells.Item(7, 8).Value2   # Producto
$ws.Cells.Item(6, 9).Value = $ws.Cells.Item(7, 9).Value2   # Categoría ID
$ws.Cells.Item(6, 10).Value = $ws.Cells.Item(7, 10).Value2 # Categoría
$ws.Cells.Item(6, 11).Value = $ws.Cells.Item(7, 11).Value2 # Variedad
$ws.Cells.Item(6, 17).Value = $ws.Cells.Item(7, 17).Value2 # Unidad de comercialización
$ws.Cells.Item(6, 18).Value = $ws.Cells.Item(7, 18).Value2 # Origen
$ws.Cells.Item(6, 20).Value = $ws.Cells.Item(7, 20).Value2 # Kg / unidad

# New values for the inserted row.
$ws.Cells.Item(6, 4).Value = 44635    # Fecha
$ws.Cells.Item(6, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 12).Value = "Segunda" # Calidad
$ws.Cells.Item(6, 13).Value = 120       # Volumen
$ws.Cells.Item(6, 14).Value = 3000      # Precio mínimo
$ws.Cells.Item(6, 15).Value = 3000      # Precio máximo
$ws.Cells.Item(6, 16).Value = 3000      # Precio promedio ponderado
$ws.Cells.Item(6, 19).Value = 1500      # Precio $/Kg
